$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "MÜŞTERİ" (customer) column in column E
# ---------------------------------------------------------------------------

# Customer value for each data row (row 2 .. row 46), in sheet order.
$customers = @(
  "HAKRO","HAKRO","COMFYBALLS","COMFYBALLS","COMFYBALLS","COMFYBALLS","COMFYBALLS",
  "COMFYBALLS","COMFYBALLS","COMFYBALLS","COMFYBALLS","COMFYBALLS","COMFYBALLS","HAKRO",
  "COMFYBALLS","COMFYBALLS","SOCIEDAD","SOCIEDAD","HAKRO","HAKRO","HAKRO","HAKRO",
  "SOCIEDAD","HAKRO","RALPH LAUREN","HAKRO","SOCIEDAD","RALPH LAUREN","RALPH LAUREN",
  "HAKRO","RALPH LAUREN","RALPH LAUREN","SELECTED FEMME","HAKRO","HAKRO","HAKRO","HAKRO",
  "HAKRO","HAKRO","SELECTED FEMME","SOCIEDAD","SOCIEDAD","HAKRO","HAKRO","HAKRO"
)

# Header cell E1: duplicate the look of the other header cells (D1) - bold
# Arial Nova 16pt on a red fill, centred - but with only a right-hand border
# (it is the new right-most column).
$d1 = $ws.Cells.Item(1, 4)
$e1 = $ws.Cells.Item(1, 5)
$d1.Copy()
$e1.PasteSpecial(-4122)
$e1.Borders.LineStyle = -4142
$e1.Borders.Item(10).LineStyle = 1
$headerText = [string][char]0x4D + [char]0xDC + [char]0x15E + [char]0x54 + [char]0x45 + [char]0x52 + [char]0x130
$e1.Value = $headerText

# Data rows
for ($i = 0; $i -lt $customers.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 5).Value = $customers[$i]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Cosmetics: resize column E to fit its contents and restore the selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).AutoFit()

$ws.Range("F33").Select()
